$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$english = @(
  'station attendant',
  'size',
  'parent',
  'final examination',
  'research',
  'garbage',
  'sugar',
  'scholarship',
  'relatives',
  'letter of recommendation',
  'graduate school',
  'typhoon',
  'day',
  'file; portfolio',
  'way; road; directions',
  'alarm clock',
  'dirty',
  'to wake (someone) up',
  'to treat (someone) to a meal',
  'to get depressed',
  'to have difficulty',
  'to take (something) out; to hand in (something)',
  'to take (someone) to (a place)',
  'to correct; to fix',
  'to become lost; to lose one''s way',
  'to be found',
  'to go to pick up',
  'to translate',
  'to laugh',
  'to iron (clothes)',
  'to collect',
  'to put (something) in',
  'to miss (a train, bus, etc.)',
  'to show',
  'to oversleep',
  'to show (someone) around',
  'to explain',
  'to come to pick up',
  'well...; let me see...',
  'by the end of the day',
  'the other day',
  'about this much',
  'I''m sorry. (casual)',
  'from now on',
  'Excuse me.; Sorry to interrupt you.',
  '(do something) by oneself',
  'in class; during the class',
  'other',
  'child',
  'to offer; to sacrifice',
  'offer',
  'the world',
  'care',
  'generation',
  'the third generation',
  'the society',
  'visibility',
  'political world',
  'limit',
  'all; entire',
  'safety',
  'whole country',
  'entirely',
  'everything; all',
  'room',
  'tennis club',
  'department manager',
  '(something) begins',
  'to begin (something)',
  'first train',
  'start',
  'every week',
  'last week',
  'one week',
  'second week',
  'weekend',
  'other than...',
  '...or more',
  '...or less',
  'within...',
  'before; formerly',
  'to think',
  'idea',
  'archeology',
  'reference',
  'to open (something)',
  '(something) opens',
  'to open',
  'opening of a store',
  'room',
  'bookstore',
  'fish shop',
  'rooftop',
  'indoor',
  'person on one''s side',
  'way of reading',
  'evening',
  'both',
  'method',
  'exercise',
  'driving',
  'lucky',
  'fate',
  'to carry',
  'to move',
  'automobile',
  'animal',
  'verb',
  'to teach',
  'classroom',
  'church',
  'Christianity',
  'textbook',
  'professor''s office',
  'basement',
  'waiting room'
)
$japanese = @(
  '駅員|えきいん（さん）',
  '大きさ|おおきさ',
  '親|おや',
  '期末試験|きまつしけん',
  '研究|けんきゅう',
  'ごみ',
  '砂糖|さとう',
  '奨学金|しょうがくきん',
  '親せき|しんせき',
  '推薦状|すいせんじょう',
  '大学院|だいがくいん',
  '台風|たいふう',
  '日|ひ',
  'ファイル',
  '道|みち',
  '目覚まし時計|めざましどけい',
  '汚い|きたない',
  '起こす|おこす',
  'おごる',
  '落ち込む|おちこむ',
  '困る|こまる',
  '出す|だす',
  '連れていく|つれていく',
  '直す|なおす',
  '道に迷う|みちにまよう',
  '見つかる|みつかる',
  '迎えに行く|むかえにいく',
  '訳す|やくす',
  '笑う|わらう',
  'アイロンをかける',
  '集める|あつめる',
  '入れる|いれる',
  '乗り遅れる|のりおくれる',
  '見せる|みせる',
  '朝寝坊する|あさねぼうする',
  '案内する|あんないする',
  '説明する|せつめいする',
  '迎えに来る|むかえにくる',
  'ええと',
  '今日中に|きょうじゅうに',
  'この間|このあいだ',
  'このぐらい',
  'ごめん',
  'これから',
  '失礼します|しつれいします',
  '自分で|じぶんで',
  '授業中に|じゅぎょうちゅうに',
  'ほかの',
  '子供|こども',
  '供える|そなえる',
  '提供|ていきょう',
  '世界|せかい',
  '世話|せわ',
  '世代|せだい',
  '三世|さんせい',
  '世の中|よのなか',
  '視界|しかい',
  '政界|せいかい',
  '限界|げんかい',
  '全部|ぜんぶ',
  '安全|あんぜん',
  '全国|ぜんこく',
  '全く|まったく',
  '全て|すべて',
  '部屋|へや',
  'テニス部|テニスぶ',
  '部長|ぶちょう',
  '始まる|はじまる',
  '始める|はじめる',
  '始発|しはつ',
  '開始|かいし',
  '毎週|まいしゅう',
  '先週|せんしゅう',
  '一週間|いっしゅうかん',
  '二週目|にしゅうめ',
  '週末|しゅうまつ',
  '～以外|～いがい',
  '～以上|～いじょう',
  '～以下|～いか',
  '～以内|～いない',
  '以前|いぜん',
  '考える|かんがえる',
  '考え|かんがえ',
  '考古学|こうこがく',
  '参考|さんこう',
  '開ける|あける',
  '開く|あく',
  '開く|ひらく',
  '開店|かいてん',
  '部屋|へや',
  '本屋|ほんや',
  '魚屋|さかなや',
  '屋上|おくじょう',
  '屋内|おくない',
  '味方|みかた',
  '読み方|よみかた',
  '夕方|ゆうがた',
  '両方|りょうほう',
  '方法|ほうほう',
  '運動|うんどう',
  '運転|うんてん',
  '運がいい|うんがいい',
  '運命|うんめい',
  '運ぶ|はこぶ',
  '動く|うごく',
  '自動車|じどうしゃ',
  '動物|どうぶつ',
  '動詞|どうし',
  '教える|おしえる',
  '教室|きょうしつ',
  '教会|きょうかい',
  'キリスト教|キリストきょう',
  '教科書|きょうかしょ',
  '研究室|けんきゅうしつ',
  '地下室|ちかしつ',
  '待合室|まちあいしつ'
)

for ($i = 0; $i -lt $english.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $english[$i]
    $ws.Cells.Item($i + 2, 2).Value = $japanese[$i]
}
